$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENTITLEMENTS")

# Update the VALIDATIONS value in J4 with the new string
$ws.Range("J4").Value = "status=200||trialSkus=DRA_TARGET_DRUG"

# Update the active selection on the sheet
$ws.Range("K2").Select()
